$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '60.791.14'
$ws.Range("E2").Value = '  -1.54%  '

# Row 3
$ws.Range("D3").Value = '3.372.75'
$ws.Range("E3").Value = '  -0.86%  '

# Row 4
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E4").Value = '  +0.16%  '

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '568.31'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -1.84%  '

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '135.87'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -1.39%  '

# Row 7
$ws.Range("E7").Value = '  -0.01%  '

# Row 8
$ws.Range("D8").Value = '3.370.65'
$ws.Range("E8").Value = '  -0.88%  '

# Row 9
$ws.Range("E9").Value = '  -1.54%  '

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '7.58'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +0.56%  '

# Row 11
$ws.Range("E11").Value = '  -3.97%  '

# Row 12
$ws.Range("E12").Value = '  -2.95%  '

# Row 13
$ws.Range("D13").Value = '3.948.38'
$ws.Range("E13").Value = '  -0.78%  '

# Row 14
$ws.Range("E14").Value = '  -0.04%  '

# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '25.88'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -0.37%  '

# Row 16
$ws.Range("D16").Value = '3.374.81'
$ws.Range("E16").Value = '  -0.61%  '

# Row 17
$ws.Range("E17").Value = '  -4.58%  '

# Row 18
$ws.Range("D18").Value = '60.911.68'
$ws.Range("E18").Value = '  -1.38%  '

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '5.78'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -1.80%  '

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '13.64'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -4.23%  '

# Row 21
$ws.Range("E21").Value = '  -2.65%  '

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '370.75'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -1.98%  '

# Row 23
$ws.Range("D23").Value = '3.507.81'
$ws.Range("E23").Value = '  -0.58%  '

# Row 24
$ws.Range("E24").Value = '  -2.54%  '

# Row 25
$ws.Range("E25").Value = '  -0.07%  '

# Row 26
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '70.68'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -0.87%  '

# Row 27
$ws.Range("E27").Value = '  -4.18%  '

# Row 28
$ws.Range("E28").Value = '  +9.53%  '

# Row 29
$ws.Range("E29").Value = '  -4.25%  '

# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -0.02%  '

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '7.28'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -3.86%  '

# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '7.99'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -3.30%  '

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '2.12'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -2.99%  '

# Row 34
$ws.Range("E34").Value = '  -0.05%  '

# Row 35
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '23.23'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -1.11%  '

# Row 36
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '5.08'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -5.04%  '

# Row 37
$ws.Range("B37").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '6.75'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -1.76%  '

# Row 38
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '1.52'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -2.57%  '

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '164.52'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -0.53%  '

# Row 40
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.0755'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -3.73%  '

# Row 41
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +0.32%  '

# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.770'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -1.58%  '

# Row 43
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '25.07'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -0.69%  '

# Row 44
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '41.75'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +0.56%  '

# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '1.70'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -2.39%  '

# Row 46
$ws.Range("E46").Value = '  -2.78%  '

# Row 47
$ws.Range("E47").Value = '  -6.10%  '

# Row 48
$ws.Range("D48").Value = '2.537.41'
$ws.Range("E48").Value = '  +8.35%  '

# Row 49
$ws.Range("E49").Value = '  +2.48%  '

# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '6.73'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -1.93%  '

# Row 51
$ws.Range("E51").Value = '  +2.13%  '
